$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 30, column D: new course entry text, using the same wrapped-text
# style already used by the other populated cells in column D (copy
# format from D2, which carries that style).
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4122) | Out-Null
$ws.Range("D30").Value = "課程 : 道場的氣場密碼-從人性到心性的實踐"

# Row 30, column B: assign the shift worker's name; the cell previously
# carried the "unfilled" yellow-highlight style, so copy the plain
# bordered style from B29 (an already-assigned name cell) before setting
# the value.
$ws.Range("B29").Copy() | Out-Null
$ws.Range("B30").PasteSpecial(-4122) | Out-Null
$ws.Range("B30").Value = "妙嘉師姐"

# Row 29, column D: short note text (keeps the existing non-wrapped style).
$ws.Range("D29").Value = "仙佛開示"

# Update the on-screen selection/scroll position to match the saved view.
$excel.ActiveWindow.ScrollRow = 24
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D31").Select() | Out-Null
